{"js": "// Load all paragraphs up front so we hold stable references to the\n// \"before\" paragraphs before any insert/delete mutates the collection.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst p = paras.items; // p[0..25], the 26 original paragraphs\n\n// --- Paragraph 0: drop the line-break + second line, keep a single line ---\np[0].insertText(\"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -23.11.24: \u26a1\ufe0f\ud83d\ude80\", Word.InsertLocation.replace);\n\n// --- Paragraph 1: new title line ---\np[1].insertText(\"Table Meets LLM: Can Large Language Models Understand Structured Table Data? A Benchmark and Empirical Study\", Word.InsertLocation.replace);\n\n// --- Paragraph 2: new intro text, followed by a new Heading2 ---\np[2].insertText(\"\u05d4\u05d9\u05d5\u05dd \u05d0\u05e0\u05d9 \u05e1\u05d5\u05e7\u05e8 \u05de\u05d0\u05de\u05e8 \u05d1\u05e0\u05d5\u05e9\u05d0 \u05e9\u05de\u05d6\u05de\u05df \u05dc\u05d0 \u05e0\u05d2\u05e2\u05ea \u05d1\u05d5(\u05d1\u05e1\u05e7\u05d9\u05e8\u05d5\u05ea) \u05d5\u05d4\u05d5\u05d0 \u05d3\u05d0\u05d8\u05d4 \u05d8\u05d1\u05dc\u05d0\u05d9. \u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05d7\u05df \u05e9\u05d0\u05dc\u05d4 \u05de\u05e8\u05ea\u05e7\u05ea - \u05d4\u05d0\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd (LLMs) \u05db\u05de\u05d5 GPT \u05d1\u05d0\u05de\u05ea \u05de\u05d1\u05d9\u05e0\u05d9\u05dd \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4 \u05d1\u05d8\u05d1\u05dc\u05d0\u05d5\u05ea?\", Word.InsertLocation.replace);\nconst h1 = p[2].insertParagraph(\"\u05e7\u05e6\u05ea: \u05e8\u05e7\u05e2\", Word.InsertLocation.after);\nh1.style = \"Heading 2\";\n\n// --- Paragraph 3: background text, followed by a new Heading2 ---\np[3].insertText(\"\u05d1\u05e9\u05e0\u05d9\u05dd \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d5\u05ea, LLMs \u05d4\u05e4\u05db\u05d5 \u05dc\u05db\u05dc\u05d9 \u05d7\u05e9\u05d5\u05d1 \u05d1\u05e2\u05d9\u05d1\u05d5\u05d3 \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea. \u05d0\u05d1\u05dc \u05d1\u05e2\u05d5\u05d3 \u05e9\u05d4\u05dd \u05de\u05e6\u05d5\u05d9\u05e0\u05d9\u05dd (\u05e1\u05d5\u05d2 \u05e9\u05dc) \u05d1\u05d4\u05d1\u05e0\u05ea \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea (\u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d8\u05e7\u05e1\u05d8), \u05d9\u05db\u05d5\u05dc\u05ea\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d9\u05d3\u05e2 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05e2\u05d3\u05d9\u05d9\u05df \u05dc\u05d0 \u05e0\u05d7\u05e7\u05e8\u05d4 \u05dc\u05e2\u05d5\u05de\u05e7 \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05de\u05d4 \u05e9\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e0\u05e1\u05d9\u05dd \u05dc\u05e2\u05e9\u05d5\u05ea \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8\", Word.InsertLocation.replace);\nconst h2 = p[3].insertParagraph(\"\u05de\u05d4 \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5?\", Word.InsertLocation.after);\nh2.style = \"Heading 2\";\n\n// --- Paragraph 4..11: the seven-task list + lead-in sentence ---\np[4].insertText(\"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e4\u05d9\u05ea\u05d7\u05d5 \u05de\u05d3\u05d3 \u05d7\u05d3\u05e9 \u05e9\u05e0\u05e7\u05e8\u05d0 (SUC (Structural Understanding Capabilities \u05e9\u05d1\u05d5\u05d7\u05df \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d1\u05e0\u05d4 \u05e9\u05dc \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea. \u05d4\u05de\u05d3\u05d3 \u05db\u05d5\u05dc\u05dc \u05e9\u05d1\u05e2 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea:\", Word.InsertLocation.replace);\np[5].insertText(\"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d8\u05d1\u05dc\u05d4\", Word.InsertLocation.replace);\np[6].insertText(\"\u05d0\u05d9\u05ea\u05d5\u05e8 \u05ea\u05d0\u05d9\u05dd \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05d9\u05dd\", Word.InsertLocation.replace);\np[7].insertText(\"\u05d7\u05d9\u05e4\u05d5\u05e9 \u05d4\u05e4\u05d5\u05da (\u05de\u05d9\u05e7\u05d5\u05dd \u05dc\u05e2\u05e8\u05da)\", Word.InsertLocation.replace);\np[8].insertText(\"\u05d0\u05d7\u05d6\u05d5\u05e8 \u05e2\u05de\u05d5\u05d3\u05d5\u05ea\", Word.InsertLocation.replace);\np[9].insertText(\"\u05d0\u05d7\u05d6\u05d5\u05e8 \u05e9\u05d5\u05e8\u05d5\u05ea\", Word.InsertLocation.replace);\np[10].insertText(\"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d5\u05d3\u05dc \u05d8\u05d1\u05dc\u05d4\", Word.InsertLocation.replace);\np[11].insertText(\"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d0\u05d9\u05dd \u05de\u05de\u05d5\u05d6\u05d2\u05d9\u05dd\", Word.InsertLocation.replace);\n\n// --- Paragraph 12: GPT-3.5/4 evaluation sentence, followed by a new Heading2 ---\np[12].insertText(\"\u05d4\u05dd \u05d1\u05d3\u05e7\u05d5 \u05d0\u05ea GPT-3.5 \u05d5-GPT-4 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d0\u05dc\u05d5 \u05ea\u05d5\u05da \u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05e4\u05d5\u05e8\u05de\u05d8\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd \u05e9\u05dc \u05e7\u05dc\u05d8 (HTML, JSON, CSV \u05d5\u05e2\u05d5\u05d3).\", Word.InsertLocation.replace);\nconst h3 = p[12].insertParagraph(\"\u05de\u05d4 \u05d4\u05dd \u05d2\u05d9\u05dc\u05d5?\", Word.InsertLocation.after);\nh3.style = \"Heading 2\";\n\n// --- Paragraph 13..16: findings ---\np[13].insertText(\"\u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05de\u05e4\u05ea\u05d9\u05e2\u05d5\u05ea! \u05d4\u05e0\u05d4 \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d4\u05e2\u05d9\u05e7\u05e8\u05d9\u05d5\u05ea:\", Word.InsertLocation.replace);\np[14].insertText(\"HTML \u05de\u05ea\u05d2\u05dc\u05d4 \u05db\u05e4\u05d5\u05e8\u05de\u05d8 \u05f4\u05d4\u05e0\u05d5\u05d7\u05f4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05dc\u05d4\u05e6\u05d2\u05ea \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05dc-LLMs\", Word.InsertLocation.replace);\np[15].insertText(\"\u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05e8\u05d0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d8\u05d5\u05d1\u05d5\u05ea \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d9\u05d7\u05e1\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea ( \u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d8\u05d1\u05dc\u05d4, \u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d0\u05d9\u05dd \u05de\u05de\u05d5\u05d6\u05d2\u05d9\u05dd) \u05d0\u05da \u05e0\u05db\u05e9\u05dc\u05d5 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea ( \u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d5\u05d3\u05dc \u05d8\u05d1\u05dc\u05d4, \u05d0\u05d7\u05d6\u05d5\u05e8 \u05e9\u05d5\u05e8\u05d4 \u05e4\u05e9\u05d5\u05d8, \u05d7\u05d9\u05e4\u05d5\u05e9 \u05ea\u05d0 \u05d1\u05d5\u05d3\u05d3)\", Word.InsertLocation.replace);\np[16].insertText(\"\u05d4\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d4\u05e9\u05ea\u05e4\u05e8\u05d5 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05e2\u05dd \u05d3\u05d5\u05d2\u05de\u05d4 \u05d0\u05d7\u05ea (one-shot) \u05dc\u05e2\u05d5\u05de\u05ea \u05d0\u05e4\u05e1 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea\", Word.InsertLocation.replace);\nconst h4 = p[16].insertParagraph(\"\u05d4\u05d7\u05d9\u05d3\u05d5\u05e9 \u05d4\u05de\u05e8\u05db\u05d6\u05d9: Self-augmented Prompting\", Word.InsertLocation.after);\nh4.style = \"Heading 2\";\n\n// --- Paragraph 17: self-augmented prompting explanation ---\np[17].insertText(\"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e4\u05d9\u05ea\u05d7\u05d5 \u05e9\u05d9\u05d8\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0\u05ea \\\"self-augmented prompting\\\" \u05e9\u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05de\u05d1\u05e7\u05e9\u05ea \u05de\u05d4\u05de\u05d5\u05d3\u05dc \u05ea\u05d7\u05d9\u05dc\u05d4 \u05dc\u05d6\u05d4\u05d5\u05ea \u05de\u05d9\u05d3\u05e2 \u05e7\u05e8\u05d9\u05d8\u05d9 \u05d1\u05d8\u05d1\u05dc\u05d4 (\u05db\u05de\u05d5 \u05d8\u05d5\u05d5\u05d7\u05d9 \u05e2\u05e8\u05db\u05d9\u05dd) \u05d5\u05d0\u05d6 \u05de\u05e9\u05ea\u05de\u05e9\u05ea \u05d1\u05de\u05d9\u05d3\u05e2 \u05d4\u05d6\u05d4 \u05db\u05d3\u05d9 \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e1\u05d5\u05e4\u05d9\u05ea. \u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d3\u05d9 \u05e8\u05e6\u05d9\u05e0\u05d9 \u05d1\u05de\u05e1\u05e4\u05e8 \u05d1\u05e0\u05e6'\u05de\u05d0\u05e8\u05e7\u05d9\u05dd)\", Word.InsertLocation.replace);\n\n// --- Paragraphs 18..22: removed entirely (old PPO conclusions + importance section) ---\np[18].delete();\np[19].delete();\np[20].delete();\np[21].delete();\np[22].delete();\n\n// --- Paragraph 23 (\"\u05e1\u05d9\u05db\u05d5\u05dd:\"): keep text, promote style to Heading2 ---\np[23].style = \"Heading 2\";\n\n// --- Paragraph 24: new closing paragraph text ---\np[24].insertText(\"\u05d0\u05e0\u05d9 \u05d7\u05d9\u05d9\u05d1 \u05dc\u05d4\u05d2\u05d9\u05d3 \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05de\u05e8\u05ea\u05e7. \u05d4\u05d5\u05d0 \u05de\u05e8\u05d0\u05d4 \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05d4\u05d4\u05ea\u05e7\u05d3\u05de\u05d5\u05ea \u05d4\u05e2\u05e6\u05d5\u05de\u05d4 \u05d1-LLMs, \u05d9\u05e9 \u05e2\u05d3\u05d9\u05d9\u05df \u05e4\u05e2\u05e8\u05d9\u05dd \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd \u05d1\u05d9\u05db\u05d5\u05dc\u05ea \u05e9\u05dc\u05d4\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4. \u05d6\u05d4 \u05de\u05d6\u05db\u05d9\u05e8 \u05dc\u05e0\u05d5 \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05de\u05e8\u05e9\u05d9\u05de\u05d9\u05dd, \u05d4\u05dd \u05e2\u05d3\u05d9\u05d9\u05df \u05e8\u05d7\u05d5\u05e7\u05d9\u05dd \u05de\u05d4\u05d1\u05e0\u05d4 \u05d0\u05e0\u05d5\u05e9\u05d9\u05ea \u05d0\u05de\u05d9\u05ea\u05d9\u05ea \u05e9\u05dc \u05de\u05d1\u05e0\u05d9\u05dd \u05d5\u05d9\u05d7\u05e1\u05d9\u05dd \u05d1\u05d9\u05df \u05d3\u05d0\u05d8\u05d4.\", Word.InsertLocation.replace);\n\n// --- Paragraph 25 (old arxiv link): replace text, then append the new closing block ---\np[25].insertText(\"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5 \u05e2\u05d1\u05d5\u05d3\u05d4 \u05dc\u05d0 \u05e8\u05e2\u05d4 \u05d1\u05e4\u05d9\u05ea\u05d5\u05d7 \u05de\u05d3\u05d3\u05d9\u05dd \u05d5\u05e9\u05d9\u05d8\u05d5\u05ea \u05e9\u05d9\u05e2\u05d6\u05e8\u05d5 \u05dc\u05e7\u05d4\u05d9\u05dc\u05d4 \u05dc\u05d4\u05de\u05e9\u05d9\u05da \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d7\u05d3\u05e9\u05d4 \u05e9\u05dc\u05d4\u05dd \u05dc-prompting \u05d4\u05d9\u05d0 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05d0\u05d1\u05dc \u05d0\u05e4\u05e7\u05d8\u05d9\u05d1\u05d9\u05ea, \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05de\u05d4 \u05e9\u05d0\u05e0\u05d7\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd - \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e4\u05e8\u05e7\u05d8\u05d9\u05d9\u05dd \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d9\u05d9\u05e9\u05dd \u05de\u05d9\u05d3.\", Word.InsertLocation.replace);\nconst h5 = p[25].insertParagraph(\"\u05de\u05d9\u05dc\u05d4 \u05d0\u05d7\u05e8\u05d5\u05e0\u05d4\", Word.InsertLocation.after);\nh5.style = \"Heading 2\";\nconst last1 = h5.insertParagraph(\"\u05d0\u05dd \u05d0\u05ea\u05dd \u05e2\u05d5\u05d1\u05d3\u05d9\u05dd \u05e2\u05dd \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05d5-LLMs, \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05d4\u05d5\u05d0 \u05d7\u05d5\u05d1\u05d4. \u05d4\u05d5\u05d0 \u05de\u05e1\u05e4\u05e7 \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05de\u05e2\u05e9\u05d9\u05d5\u05ea \u05d5\u05db\u05dc\u05d9\u05dd \u05e9\u05d9\u05de\u05d5\u05e9\u05d9\u05d9\u05dd. \u05d4\u05e7\u05d5\u05d3 \u05d5\u05d4\u05d3\u05d0\u05d8\u05d4 \u05d6\u05de\u05d9\u05e0\u05d9\u05dd \u05d1-GitHub, \u05d0\u05d6 \u05d0\u05ea\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d4\u05ea\u05d7\u05d9\u05dc \u05dc\u05e9\u05d7\u05e7 \u05e2\u05dd \u05d6\u05d4 \u05d9\u05e9\u05e8.\", Word.InsertLocation.after);\nlast1.style = \"Normal\";\nconst last2 = last1.insertParagraph(\"\u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05d9\u05d4\u05d9\u05d4 \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05d9\u05da \u05d4\u05de\u05de\u05e6\u05d0\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05d9\u05e9\u05e4\u05d9\u05e2\u05d5 \u05e2\u05dc \u05d4\u05d3\u05d5\u05e8 \u05d4\u05d1\u05d0 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d4\u05d0\u05dd \u05e0\u05e8\u05d0\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05de\u05ea\u05d5\u05db\u05e0\u05e0\u05d9\u05dd \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05dc\u05d4\u05d1\u05e0\u05ea \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4?\", Word.InsertLocation.after);\nlast2.style = \"Normal\";\nconst last3 = last2.insertParagraph(\"https://arxiv.org/abs/2305.13062\", Word.InsertLocation.after);\nlast3.style = \"Normal\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Process paragraphs from the LAST original index down to the FIRST so that\n# inserting/deleting paragraphs never shifts the index of a paragraph we\n# still have left to touch (the original document has exactly 26 paragraphs).\n\n# --- Paragraph 26 (old arxiv link) -> new closing sentence, then 4 new paragraphs appended ---\n$p26 = $d.Paragraphs.Item(26)\n$p26.Range.Text = \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5 \u05e2\u05d1\u05d5\u05d3\u05d4 \u05dc\u05d0 \u05e8\u05e2\u05d4 \u05d1\u05e4\u05d9\u05ea\u05d5\u05d7 \u05de\u05d3\u05d3\u05d9\u05dd \u05d5\u05e9\u05d9\u05d8\u05d5\u05ea \u05e9\u05d9\u05e2\u05d6\u05e8\u05d5 \u05dc\u05e7\u05d4\u05d9\u05dc\u05d4 \u05dc\u05d4\u05de\u05e9\u05d9\u05da \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d4. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d4\u05d7\u05d3\u05e9\u05d4 \u05e9\u05dc\u05d4\u05dd \u05dc-prompting \u05d4\u05d9\u05d0 \u05e4\u05e9\u05d5\u05d8\u05d4 \u05d0\u05d1\u05dc \u05d0\u05e4\u05e7\u05d8\u05d9\u05d1\u05d9\u05ea, \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05de\u05d4 \u05e9\u05d0\u05e0\u05d7\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd - \u05e4\u05ea\u05e8\u05d5\u05e0\u05d5\u05ea \u05e4\u05e8\u05e7\u05d8\u05d9\u05d9\u05dd \u05e9\u05d0\u05e4\u05e9\u05e8 \u05dc\u05d9\u05d9\u05e9\u05dd \u05de\u05d9\u05d3.\"\n$p26.Range.InsertParagraphAfter()\n$p27 = $d.Paragraphs.Item(27)\n$p27.Range.Text = \"\u05de\u05d9\u05dc\u05d4 \u05d0\u05d7\u05e8\u05d5\u05e0\u05d4\"\n$p27.Range.Style = \"Heading 2\"\n$p27.Range.InsertParagraphAfter()\n$p28 = $d.Paragraphs.Item(28)\n$p28.Range.Text = \"\u05d0\u05dd \u05d0\u05ea\u05dd \u05e2\u05d5\u05d1\u05d3\u05d9\u05dd \u05e2\u05dd \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05d5-LLMs, \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05d4\u05d5\u05d0 \u05d7\u05d5\u05d1\u05d4. \u05d4\u05d5\u05d0 \u05de\u05e1\u05e4\u05e7 \u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05de\u05e2\u05e9\u05d9\u05d5\u05ea \u05d5\u05db\u05dc\u05d9\u05dd \u05e9\u05d9\u05de\u05d5\u05e9\u05d9\u05d9\u05dd. \u05d4\u05e7\u05d5\u05d3 \u05d5\u05d4\u05d3\u05d0\u05d8\u05d4 \u05d6\u05de\u05d9\u05e0\u05d9\u05dd \u05d1-GitHub, \u05d0\u05d6 \u05d0\u05ea\u05dd \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d4\u05ea\u05d7\u05d9\u05dc \u05dc\u05e9\u05d7\u05e7 \u05e2\u05dd \u05d6\u05d4 \u05d9\u05e9\u05e8.\"\n$p28.Range.Style = \"Normal\"\n$p28.Range.InsertParagraphAfter()\n$p29 = $d.Paragraphs.Item(29)\n$p29.Range.Text = \"\u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05d9\u05d4\u05d9\u05d4 \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05d9\u05da \u05d4\u05de\u05de\u05e6\u05d0\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05d9\u05e9\u05e4\u05d9\u05e2\u05d5 \u05e2\u05dc \u05d4\u05d3\u05d5\u05e8 \u05d4\u05d1\u05d0 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05d4\u05d0\u05dd \u05e0\u05e8\u05d0\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05de\u05ea\u05d5\u05db\u05e0\u05e0\u05d9\u05dd \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05dc\u05d4\u05d1\u05e0\u05ea \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4?\"\n$p29.Range.Style = \"Normal\"\n$p29.Range.InsertParagraphAfter()\n$p30 = $d.Paragraphs.Item(30)\n$p30.Range.Text = \"https://arxiv.org/abs/2305.13062\"\n$p30.Range.Style = \"Normal\"\n\n# --- Paragraph 25 -> new closing paragraph text ---\n$p25 = $d.Paragraphs.Item(25)\n$p25.Range.Text = \"\u05d0\u05e0\u05d9 \u05d7\u05d9\u05d9\u05d1 \u05dc\u05d4\u05d2\u05d9\u05d3 \u05e9\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05de\u05e8\u05ea\u05e7. \u05d4\u05d5\u05d0 \u05de\u05e8\u05d0\u05d4 \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05d4\u05d4\u05ea\u05e7\u05d3\u05de\u05d5\u05ea \u05d4\u05e2\u05e6\u05d5\u05de\u05d4 \u05d1-LLMs, \u05d9\u05e9 \u05e2\u05d3\u05d9\u05d9\u05df \u05e4\u05e2\u05e8\u05d9\u05dd \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd \u05d1\u05d9\u05db\u05d5\u05dc\u05ea \u05e9\u05dc\u05d4\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4. \u05d6\u05d4 \u05de\u05d6\u05db\u05d9\u05e8 \u05dc\u05e0\u05d5 \u05e9\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05d0\u05dc\u05d4 \u05de\u05e8\u05e9\u05d9\u05de\u05d9\u05dd, \u05d4\u05dd \u05e2\u05d3\u05d9\u05d9\u05df \u05e8\u05d7\u05d5\u05e7\u05d9\u05dd \u05de\u05d4\u05d1\u05e0\u05d4 \u05d0\u05e0\u05d5\u05e9\u05d9\u05ea \u05d0\u05de\u05d9\u05ea\u05d9\u05ea \u05e9\u05dc \u05de\u05d1\u05e0\u05d9\u05dd \u05d5\u05d9\u05d7\u05e1\u05d9\u05dd \u05d1\u05d9\u05df \u05d3\u05d0\u05d8\u05d4.\"\n\n# --- Paragraph 24 (\"\u05e1\u05d9\u05db\u05d5\u05dd:\") -> keep text, promote style to Heading2 ---\n$p24 = $d.Paragraphs.Item(24)\n$p24.Range.Style = \"Heading 2\"\n\n# --- Paragraphs 19-23 -> removed entirely ---\n$d.Paragraphs.Item(23).Range.Delete()\n$d.Paragraphs.Item(22).Range.Delete()\n$d.Paragraphs.Item(21).Range.Delete()\n$d.Paragraphs.Item(20).Range.Delete()\n$d.Paragraphs.Item(19).Range.Delete()\n\n# --- Paragraph 18 -> self-augmented prompting explanation ---\n$p18 = $d.Paragraphs.Item(18)\n$p18.Range.Text = \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e4\u05d9\u05ea\u05d7\u05d5 \u05e9\u05d9\u05d8\u05d4 \u05d7\u05d3\u05e9\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0\u05ea \"\"self-augmented prompting\"\" \u05e9\u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05de\u05d1\u05e7\u05e9\u05ea \u05de\u05d4\u05de\u05d5\u05d3\u05dc \u05ea\u05d7\u05d9\u05dc\u05d4 \u05dc\u05d6\u05d4\u05d5\u05ea \u05de\u05d9\u05d3\u05e2 \u05e7\u05e8\u05d9\u05d8\u05d9 \u05d1\u05d8\u05d1\u05dc\u05d4 (\u05db\u05de\u05d5 \u05d8\u05d5\u05d5\u05d7\u05d9 \u05e2\u05e8\u05db\u05d9\u05dd) \u05d5\u05d0\u05d6 \u05de\u05e9\u05ea\u05de\u05e9\u05ea \u05d1\u05de\u05d9\u05d3\u05e2 \u05d4\u05d6\u05d4 \u05db\u05d3\u05d9 \u05dc\u05e9\u05e4\u05e8 \u05d0\u05ea \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e1\u05d5\u05e4\u05d9\u05ea. \u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05e9\u05d9\u05e4\u05d5\u05e8 \u05d3\u05d9 \u05e8\u05e6\u05d9\u05e0\u05d9 \u05d1\u05de\u05e1\u05e4\u05e8 \u05d1\u05e0\u05e6'\u05de\u05d0\u05e8\u05e7\u05d9\u05dd)\"\n\n# --- Paragraph 17 -> one-shot sentence, then new Heading2 after ---\n$p17 = $d.Paragraphs.Item(17)\n$p17.Range.Text = \"\u05d4\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d4\u05e9\u05ea\u05e4\u05e8\u05d5 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea \u05e2\u05dd \u05d3\u05d5\u05d2\u05de\u05d4 \u05d0\u05d7\u05ea (one-shot) \u05dc\u05e2\u05d5\u05de\u05ea \u05d0\u05e4\u05e1 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea\"\n$p17.Range.InsertParagraphAfter()\n$p17b = $d.Paragraphs.Item(18)\n$p17b.Range.Text = \"\u05d4\u05d7\u05d9\u05d3\u05d5\u05e9 \u05d4\u05de\u05e8\u05db\u05d6\u05d9: Self-augmented Prompting\"\n$p17b.Range.Style = \"Heading 2\"\n\n# --- Paragraph 16 -> complex/simple tasks sentence ---\n$p16 = $d.Paragraphs.Item(16)\n$p16.Range.Text = \"\u05d4\u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d4\u05e8\u05d0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05d8\u05d5\u05d1\u05d5\u05ea \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d9\u05d7\u05e1\u05d9\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea ( \u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d8\u05d1\u05dc\u05d4, \u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d0\u05d9\u05dd \u05de\u05de\u05d5\u05d6\u05d2\u05d9\u05dd) \u05d0\u05da \u05e0\u05db\u05e9\u05dc\u05d5 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea ( \u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d5\u05d3\u05dc \u05d8\u05d1\u05dc\u05d4, \u05d0\u05d7\u05d6\u05d5\u05e8 \u05e9\u05d5\u05e8\u05d4 \u05e4\u05e9\u05d5\u05d8, \u05d7\u05d9\u05e4\u05d5\u05e9 \u05ea\u05d0 \u05d1\u05d5\u05d3\u05d3)\"\n\n# --- Paragraph 15 -> HTML format sentence ---\n$p15 = $d.Paragraphs.Item(15)\n$p15.Range.Text = \"HTML \u05de\u05ea\u05d2\u05dc\u05d4 \u05db\u05e4\u05d5\u05e8\u05de\u05d8 \u05f4\u05d4\u05e0\u05d5\u05d7\u05f4 \u05d1\u05d9\u05d5\u05ea\u05e8 \u05dc\u05d4\u05e6\u05d2\u05ea \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05dc-LLMs\"\n\n# --- Paragraph 14 -> results intro sentence ---\n$p14 = $d.Paragraphs.Item(14)\n$p14.Range.Text = \"\u05d4\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05de\u05e4\u05ea\u05d9\u05e2\u05d5\u05ea! \u05d4\u05e0\u05d4 \u05d4\u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d4\u05e2\u05d9\u05e7\u05e8\u05d9\u05d5\u05ea:\"\n\n# --- Paragraph 13 -> GPT sentence, then new Heading2 after ---\n$p13 = $d.Paragraphs.Item(13)\n$p13.Range.Text = \"\u05d4\u05dd \u05d1\u05d3\u05e7\u05d5 \u05d0\u05ea GPT-3.5 \u05d5-GPT-4 \u05d1\u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05d0\u05dc\u05d5 \u05ea\u05d5\u05da \u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05e4\u05d5\u05e8\u05de\u05d8\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd \u05e9\u05dc \u05e7\u05dc\u05d8 (HTML, JSON, CSV \u05d5\u05e2\u05d5\u05d3).\"\n$p13.Range.InsertParagraphAfter()\n$p13b = $d.Paragraphs.Item(14)\n$p13b.Range.Text = \"\u05de\u05d4 \u05d4\u05dd \u05d2\u05d9\u05dc\u05d5?\"\n$p13b.Range.Style = \"Heading 2\"\n\n# --- Paragraphs 6-12 -> seven-task list ---\n$d.Paragraphs.Item(12).Range.Text = \"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05ea\u05d0\u05d9\u05dd \u05de\u05de\u05d5\u05d6\u05d2\u05d9\u05dd\"\n$d.Paragraphs.Item(11).Range.Text = \"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d5\u05d3\u05dc \u05d8\u05d1\u05dc\u05d4\"\n$d.Paragraphs.Item(10).Range.Text = \"\u05d0\u05d7\u05d6\u05d5\u05e8 \u05e9\u05d5\u05e8\u05d5\u05ea\"\n$d.Paragraphs.Item(9).Range.Text = \"\u05d0\u05d7\u05d6\u05d5\u05e8 \u05e2\u05de\u05d5\u05d3\u05d5\u05ea\"\n$d.Paragraphs.Item(8).Range.Text = \"\u05d7\u05d9\u05e4\u05d5\u05e9 \u05d4\u05e4\u05d5\u05da (\u05de\u05d9\u05e7\u05d5\u05dd \u05dc\u05e2\u05e8\u05da)\"\n$d.Paragraphs.Item(7).Range.Text = \"\u05d0\u05d9\u05ea\u05d5\u05e8 \u05ea\u05d0\u05d9\u05dd \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05d9\u05dd\"\n$d.Paragraphs.Item(6).Range.Text = \"\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d2\u05d1\u05d5\u05dc\u05d5\u05ea \u05d8\u05d1\u05dc\u05d4\"\n\n# --- Paragraph 5 -> SUC metric sentence ---\n$p5 = $d.Paragraphs.Item(5)\n$p5.Range.Text = \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e4\u05d9\u05ea\u05d7\u05d5 \u05de\u05d3\u05d3 \u05d7\u05d3\u05e9 \u05e9\u05e0\u05e7\u05e8\u05d0 (SUC (Structural Understanding Capabilities \u05e9\u05d1\u05d5\u05d7\u05df \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d1\u05e0\u05d4 \u05e9\u05dc \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea. \u05d4\u05de\u05d3\u05d3 \u05db\u05d5\u05dc\u05dc \u05e9\u05d1\u05e2 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea:\"\n\n# --- Paragraph 4 -> background sentence, then new Heading2 after ---\n$p4 = $d.Paragraphs.Item(4)\n$p4.Range.Text = \"\u05d1\u05e9\u05e0\u05d9\u05dd \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d5\u05ea, LLMs \u05d4\u05e4\u05db\u05d5 \u05dc\u05db\u05dc\u05d9 \u05d7\u05e9\u05d5\u05d1 \u05d1\u05e2\u05d9\u05d1\u05d5\u05d3 \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea. \u05d0\u05d1\u05dc \u05d1\u05e2\u05d5\u05d3 \u05e9\u05d4\u05dd \u05de\u05e6\u05d5\u05d9\u05e0\u05d9\u05dd (\u05e1\u05d5\u05d2 \u05e9\u05dc) \u05d1\u05d4\u05d1\u05e0\u05ea \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea (\u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d8\u05e7\u05e1\u05d8), \u05d9\u05db\u05d5\u05dc\u05ea\u05dd \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05d9\u05d3\u05e2 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d8\u05d1\u05dc\u05d0\u05d5\u05ea \u05e2\u05d3\u05d9\u05d9\u05df \u05dc\u05d0 \u05e0\u05d7\u05e7\u05e8\u05d4 \u05dc\u05e2\u05d5\u05de\u05e7 \u05d5\u05d6\u05d4 \u05d1\u05d3\u05d9\u05d5\u05e7 \u05de\u05d4 \u05e9\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e0\u05e1\u05d9\u05dd \u05dc\u05e2\u05e9\u05d5\u05ea \u05d1\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8\"\n$p4.Range.InsertParagraphAfter()\n$p4b = $d.Paragraphs.Item(5)\n$p4b.Range.Text = \"\u05de\u05d4 \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05e2\u05e9\u05d5?\"\n$p4b.Range.Style = \"Heading 2\"\n\n# --- Paragraph 3 -> intro sentence, then new Heading2 after ---\n$p3 = $d.Paragraphs.Item(3)\n$p3.Range.Text = \"\u05d4\u05d9\u05d5\u05dd \u05d0\u05e0\u05d9 \u05e1\u05d5\u05e7\u05e8 \u05de\u05d0\u05de\u05e8 \u05d1\u05e0\u05d5\u05e9\u05d0 \u05e9\u05de\u05d6\u05de\u05df \u05dc\u05d0 \u05e0\u05d2\u05e2\u05ea \u05d1\u05d5(\u05d1\u05e1\u05e7\u05d9\u05e8\u05d5\u05ea) \u05d5\u05d4\u05d5\u05d0 \u05d3\u05d0\u05d8\u05d4 \u05d8\u05d1\u05dc\u05d0\u05d9. \u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05d7\u05df \u05e9\u05d0\u05dc\u05d4 \u05de\u05e8\u05ea\u05e7\u05ea - \u05d4\u05d0\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd (LLMs) \u05db\u05de\u05d5 GPT \u05d1\u05d0\u05de\u05ea \u05de\u05d1\u05d9\u05e0\u05d9\u05dd \u05de\u05d9\u05d3\u05e2 \u05de\u05d5\u05d1\u05e0\u05d4 \u05d1\u05d8\u05d1\u05dc\u05d0\u05d5\u05ea?\"\n$p3.Range.InsertParagraphAfter()\n$p3b = $d.Paragraphs.Item(4)\n$p3b.Range.Text = \"\u05e7\u05e6\u05ea: \u05e8\u05e7\u05e2\"\n$p3b.Range.Style = \"Heading 2\"\n\n# --- Paragraph 2 -> new title line ---\n$p2 = $d.Paragraphs.Item(2)\n$p2.Range.Text = \"Table Meets LLM: Can Large Language Models Understand Structured Table Data? A Benchmark and Empirical Study\"\n\n# --- Paragraph 1 -> drop line-break + second line, single line only ---\n$p1 = $d.Paragraphs.Item(1)\n$p1.Range.Text = \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -23.11.24: \u26a1\ufe0f\ud83d\ude80\"\n"}
